$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valores colhidos")

# Update raw input values
$ws.Range("A2").Value = 222
$ws.Range("D2").Value = 212
$ws.Range("C2").Formula = "=B2/1000"

$ws.Range("D3").Value = 288
$ws.Range("D4").Value = 354
$ws.Range("D5").Value = 405

# Recalculate all formulas so dependent cells (B3, B4, B5, E2:E5, F2:F5, H2:H5, I2:I5, F7) update
$excel.Calculate()

# Update selection to match the committed state
$ws.Activate()
$ws.Range("F9").Select()
